$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Status text change for the "66e31fa5..." row (row 3 in each sheet):
# "Ready for handoff" -> "Handback transform failed"
$overview.Range("B3").Value = "Handback transform failed"
$overview.Range("C3").Value = "Handback transform failed"
$zhcn.Range("C3").Value = "Handback transform failed"
$dede.Range("C3").Value = "Handback transform failed"

# New Error Detail (column K) values for row 3 on the locale sheets
$zhcn.Range("K3").Value = "Handback file name: pp4jnmga.5uo is different with handoff file name: 66e31fa5-54e9-4ce6-b319-1b1796b036a3.c14b791b01928c9ba11ba9c92dc002792c00f4b0.zh-cn."
$dede.Range("K3").Value = "Handback file name: pp4jnmga.5uo is different with handoff file name: 66e31fa5-54e9-4ce6-b319-1b1796b036a3.c14b791b01928c9ba11ba9c92dc002792c00f4b0.de-de."
